# Apply the "Trade #19 closed" update to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200          # Current Capital
$summary.Range("B4").Value = 0             # Total P&L $
$summary.Range("B5").Value = 0             # Total P&L %
$summary.Range("B6").Value = 19            # Total Trades
$summary.Range("B8").Value = 8             # Losing Trades
$summary.Range("B9").Value = 31.58         # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100            # Capital
$status.Range("D4").Value = 19             # Trades
$status.Range("E4").Value = 0              # P&L $
$status.Range("F4").Value = -0             # P&L %
$status.Range("G4").Value = 31.58          # Win Rate %

# ---------------------------------------------------------------
# Helper that appends the closed Trade #19 row to a trades sheet
# ---------------------------------------------------------------
function Add-Trade19Row($ws) {
    $ws.Cells.Item(20, 1).Value = 19
    $ws.Cells.Item(20, 2).Value = "'2026-02-17"
    $ws.Cells.Item(20, 3).Value = "'04:08:00"
    $ws.Cells.Item(20, 4).Value = "MarketMaking"
    $ws.Cells.Item(20, 5).Value = "DOWN"
    $ws.Cells.Item(20, 6).Value = 0.73
    $ws.Cells.Item(20, 7).Value = 0.65
    $ws.Cells.Item(20, 8).Value = "CLOSED"
    $ws.Cells.Item(20, 9).Value = -10.9589
    $ws.Cells.Item(20, 10).Value = -0.08
    $ws.Cells.Item(20, 11).Value = 100
    $ws.Cells.Item(20, 12).Value = 0
    $ws.Cells.Item(20, 13).Value = 0
    $ws.Cells.Item(20, 14).Value = 0.6
    $ws.Cells.Item(20, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(20, 16).Value = "early_exit"
    $ws.Cells.Item(20, 17).Value = 0.12
}

# ---------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade19Row $allTrades

# ---------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade19Row $marketMaking
